# Add a new "time_taken" metadata column (F) to the panel data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1, using the same header formatting (bold/border/centered)
# already applied to the other header cells (e.g. E1), by copying that
# cell's formatting instead of constructing a brand-new style.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Per-row "time taken" timestamps for the data rows.
$ws.Range("F2").Value = "2021-10-05 10:50:59.566891"
$ws.Range("F3").Value = "2021-10-05 10:50:59.566904"
$ws.Range("F4").Value = "2021-10-05 10:50:59.566909"
$ws.Range("F5").Value = "2021-10-05 10:50:59.566912"
$ws.Range("F6").Value = "2021-10-05 10:50:59.566916"
$ws.Range("F7").Value = "2021-10-05 10:50:59.566919"
$ws.Range("F8").Value = "2021-10-05 10:50:59.566922"
